$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 1325.625
$ws.Range("J32").Value = 1150.6666
$ws.Range("L32").Value = 1150.6666
$ws.Range("N32").Value = -1802.6666
# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 16597.715
$ws.Range("I43").Value = 100000
$ws.Range("J43").Value = 2697.3333
$ws.Range("K43").Value = 100000
$ws.Range("L43").Value = 2697.3333
$ws.Range("M43").Value = -99931
$ws.Range("N43").Value = -2835.3333
# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 19440
$ws.Range("J51").Value = 19440
$ws.Range("L51").Value = 19440
$ws.Range("N51").Value = -20408
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 22328.445
$ws.Range("I64").Value = 3786.3333
$ws.Range("K64").Value = 3786.3333
$ws.Range("M64").Value = -3538.3333
# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 22328.445
$ws.Range("I67").Value = 3786.3333
$ws.Range("K67").Value = 3786.3333
$ws.Range("M67").Value = -2928.3333
# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 3280.2
$ws.Range("I70").Value = 3533.3333
$ws.Range("J70").Value = 3171.7144
$ws.Range("K70").Value = 10599.9999
$ws.Range("L70").Value = 9515.143199999999
$ws.Range("M70").Value = -10329.9999
$ws.Range("N70").Value = -10055.1432
# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 3280.2
$ws.Range("I73").Value = 3533.3333
$ws.Range("J73").Value = 3171.7144
$ws.Range("K73").Value = 10599.9999
$ws.Range("L73").Value = 9515.143199999999
$ws.Range("M73").Value = -9663.999899999999
$ws.Range("N73").Value = -11387.1432
# Row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 2781393.8
$ws.Range("I74").Value = 3229376.2
$ws.Range("K74").Value = 3229376.2
$ws.Range("M74").Value = -3228440.2
# Row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 2781393.8
$ws.Range("I77").Value = 3229376.2
$ws.Range("K77").Value = 16146881
$ws.Range("M77").Value = -16142201
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 4053.2646
$ws.Range("I137").Value = 925.8
$ws.Range("J137").Value = 12740.667
$ws.Range("K137").Value = 2777.4
$ws.Range("L137").Value = 38222.001
$ws.Range("M137").Value = -227.3999999999996
$ws.Range("N137").Value = -43322.001

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1131.4286
$ws.Range("I2").Value = 903.3333
$ws.Range("K2").Value = 903.3333
$ws.Range("M2").Value = -790.3333
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 160890.6
$ws.Range("I32").Value = 158758.05
$ws.Range("J32").Value = 183637.83
$ws.Range("K32").Value = 158758.05
$ws.Range("L32").Value = 183637.83
$ws.Range("M32").Value = -158471.05
$ws.Range("N32").Value = -184211.83
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1471.6444
$ws.Range("I45").Value = 1269.3334
$ws.Range("J45").Value = 2280.889
$ws.Range("K45").Value = 1269.3334
$ws.Range("L45").Value = 2280.889
$ws.Range("M45").Value = -892.3334
$ws.Range("N45").Value = -3034.889
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1010.9091
$ws.Range("I110").Value = 731.7143
$ws.Range("J110").Value = 1499.5
$ws.Range("K110").Value = 731.7143
$ws.Range("L110").Value = 1499.5
$ws.Range("M110").Value = 1313.2857
$ws.Range("N110").Value = -5589.5
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1131.4286
$ws.Range("I116").Value = 903.3333
$ws.Range("K116").Value = 903.3333
$ws.Range("M116").Value = 1390.6667

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1131.4286
$ws.Range("I3").Value = 903.3333
$ws.Range("K3").Value = 903.3333
$ws.Range("M3").Value = -789.3333
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 1380
$ws.Range("I105").Value = 1433.3334
$ws.Range("K105").Value = 1433.3334
$ws.Range("M105").Value = 313.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 17387.11
$ws.Range("I31").Value = 19248.186
$ws.Range("J31").Value = 10271.235
$ws.Range("K31").Value = 19248.186
$ws.Range("L31").Value = 10271.235
$ws.Range("M31").Value = -18953.186
$ws.Range("N31").Value = -10861.235
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 17387.11
$ws.Range("I34").Value = 19248.186
$ws.Range("J34").Value = 10271.235
$ws.Range("K34").Value = 19248.186
$ws.Range("L34").Value = 10271.235
$ws.Range("M34").Value = -19046.186
$ws.Range("N34").Value = -10675.235
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 925
$ws.Range("I105").Value = 940
$ws.Range("J105").Value = 895
$ws.Range("K105").Value = 940
$ws.Range("L105").Value = 895
$ws.Range("M105").Value = 807
$ws.Range("N105").Value = -4389

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 1018.2
$ws.Range("I2").Value = 18.857143
$ws.Range("J2").Value = 3350
$ws.Range("K2").Value = 18.857143
$ws.Range("L2").Value = 3350
$ws.Range("M2").Value = 94.14285699999999
$ws.Range("N2").Value = -3576
# Row 51 (Leve Item ID 27222)
$ws.Range("H51").Value = 30884
$ws.Range("J51").Value = 30884
$ws.Range("L51").Value = 30884
$ws.Range("N51").Value = -31902
# Row 95 (Leve Item ID 18235)
$ws.Range("H95").Value = 12000
$ws.Range("J95").Value = 12000
$ws.Range("L95").Value = 12000
$ws.Range("N95").Value = -17492
# Row 96 (Leve Item ID 18261)
$ws.Range("H96").Value = 9000
$ws.Range("J96").Value = 9000
$ws.Range("L96").Value = 9000
$ws.Range("N96").Value = -14492
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 1059.1034
$ws.Range("I97").Value = 990.9474
$ws.Range("K97").Value = 990.9474
$ws.Range("M97").Value = -494.9474
# Row 98 (Leve Item ID 18359)
$ws.Range("H98").Value = 43514.332
$ws.Range("J98").Value = 43514.332
$ws.Range("L98").Value = 43514.332
$ws.Range("N98").Value = -49504.332
# Row 99 (Leve Item ID 19532)
$ws.Range("H99").Value = 9500
$ws.Range("I99").Value = 4375
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 4375
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -2129
$ws.Range("N99").Value = -34492
# Row 110 (Leve Item ID 25802)
$ws.Range("H110").Value = 37851
$ws.Range("J110").Value = 37851
$ws.Range("L110").Value = 37851
$ws.Range("N110").Value = -46031

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 818
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 896.6667
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 896.6667
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1272.6667
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 1639
$ws.Range("I93").Value = 1125
$ws.Range("J93").Value = 1810.3334
$ws.Range("K93").Value = 1125
$ws.Range("L93").Value = 1810.3334
$ws.Range("M93").Value = 123
$ws.Range("N93").Value = -4306.3334
